# [ADDITIONAL SCRAPING] Add a "Player Info" sheet and replace the scraped
# match-card URLs with bare match codes on the batting/bowling sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "Player Info" sheet in front of the existing sheets.
# ---------------------------------------------------------------------
$battingSheetTmp = $wb.Worksheets.Item("ODI Batting")
$playerSheet = $wb.Worksheets.Add()
$playerSheet.Move($battingSheetTmp)
$playerSheet.Name = "Player Info"

# Re-fetch the other sheets by name now that the sheet collection has
# shifted (index-based handles taken before the Add/Move are stale).
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $playerSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$headerRange = $playerSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Write the ID as text (not a number) to mirror the scraped inline-string
# source data, then drop back to the default style so only the value's
# type changes, not the cell's formatting.
$playerSheet.Range("A2").NumberFormat = "@"
$playerSheet.Range("A2").Value = "4694"
$playerSheet.Range("A2").Style = "Normal"

$playerSheet.Range("B2").Value = "Dawid Johannes Malan"
$playerSheet.Range("C2").Value = "Left Handed"
$playerSheet.Range("D2").Value = "Right Arm Leg Break"

# ---------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code.
# ---------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingLastRow = $battingSheet.Cells.Item($battingSheet.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = $cell.Text
    if ($link -and $link -match "MatchCode=(.+)$") {
        $cell.NumberFormat = "@"
        $cell.Value = $matches[1]
        $cell.Style = "Normal"
    }
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code.
# ---------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingLastRow = $bowlingSheet.Cells.Item($bowlingSheet.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $link = $cell.Text
    if ($link -and $link -match "MatchCode=(.+)$") {
        $cell.NumberFormat = "@"
        $cell.Value = $matches[1]
        $cell.Style = "Normal"
    }
}

Write-Host "Edit complete"
